$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:E51 as Text first so numeric-looking strings (e.g. "1.00", "81.393.02")
# are preserved verbatim as text instead of being auto-coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '81.393.02'
$ws.Range("E2").Value = '  +5.24%  '

# Row 3
$ws.Range("D3").Value = '3.203.48'
$ws.Range("E3").Value = '  +1.88%  '

# Row 4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '211.14'
$ws.Range("E5").Value = '  +3.22%  '

# Row 6
$ws.Range("D6").Value = '638.77'
$ws.Range("E6").Value = '  +1.47%  '

# Row 7
$ws.Range("D7").Value = '0.295'
$ws.Range("E7").Value = '  +30.68%  '

# Row 9
$ws.Range("D9").Value = '0.596'
$ws.Range("E9").Value = '  +4.31%  '

# Row 10
$ws.Range("D10").Value = '3.200.23'
$ws.Range("E10").Value = '  +1.83%  '

# Row 11
$ws.Range("D11").Value = '0.599'
$ws.Range("E11").Value = '  +15.53%  '

# Row 12
$ws.Range("D12").Value = '0.0000269'
$ws.Range("E12").Value = '  +20.33%  '

# Row 13
$ws.Range("D13").Value = '0.166'
$ws.Range("E13").Value = '  +2.37%  '

# Row 14
$ws.Range("D14").Value = '5.43'
$ws.Range("E14").Value = '  +0.81%  '

# Row 15
$ws.Range("D15").Value = '3.794.47'
$ws.Range("E15").Value = '  +2.03%  '

# Row 16
$ws.Range("D16").Value = '32.46'
$ws.Range("E16").Value = '  +6.00%  '

# Row 17
$ws.Range("D17").Value = '81.285.87'
$ws.Range("E17").Value = '  +5.33%  '

# Row 18
$ws.Range("D18").Value = '3.200.07'
$ws.Range("E18").Value = '  +1.84%  '

# Row 19
$ws.Range("D19").Value = '14.53'
$ws.Range("E19").Value = '  +4.07%  '

# Row 20
$ws.Range("D20").Value = '3.18'
$ws.Range("E20").Value = '  +11.53%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '9.34'
$ws.Range("E21").Value = '  +1.38%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '445.35'
$ws.Range("E22").Value = '  +7.98%  '

# Row 23
$ws.Range("D23").Value = '5.30'
$ws.Range("E23").Value = '  +12.96%  '

# Row 24
$ws.Range("D24").Value = '7.10'
$ws.Range("E24").Value = '  +6.70%  '

# Row 25
$ws.Range("D25").Value = '5.12'
$ws.Range("E25").Value = '  +10.65%  '

# Row 26
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '11.35'
$ws.Range("E26").Value = '  +8.29%  '

# Row 27
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '3.367.13'
$ws.Range("E27").Value = '  +1.90%  '

# Row 28
$ws.Range("B28").Value = 'Litecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D28").Value = '77.57'
$ws.Range("E28").Value = '  +3.68%  '

# Row 29
$ws.Range("D29").Value = '0.0000129'
$ws.Range("E29").Value = '  +12.81%  '

# Row 30
$ws.Range("E30").Value = '  -0.16%  '

# Row 31
$ws.Range("D31").Value = '9.27'
$ws.Range("E31").Value = '  +6.12%  '

# Row 32
$ws.Range("E32").Value = '  +0.48%  '

# Row 33
$ws.Range("D33").Value = '571.67'
$ws.Range("E33").Value = '  +9.78%  '

# Row 34
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").Value = '  +4.66%  '

# Row 35
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '0.154'
$ws.Range("E35").Value = '  +13.92%  '

# Row 36
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").Value = '2.05'
$ws.Range("E36").Value = '  +5.13%  '

# Row 37
$ws.Range("B37").Value = 'Cronos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D37").Value = '0.130'
$ws.Range("E37").Value = '  +23.64%  '

# Row 38
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '23.31'
$ws.Range("E38").Value = '  +7.35%  '

# Row 39
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").Value = '0.418'
$ws.Range("E39").Value = '  +6.45%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.16%  '

# Row 41
$ws.Range("D41").Value = '2.11'
$ws.Range("E41").Value = '  +21.18%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '6.07'
$ws.Range("E42").Value = '  +13.53%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '3.12'
$ws.Range("E43").Value = '  +23.27%  '

# Row 44
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").Value = '20.81'
$ws.Range("E44").Value = '  +3.69%  '

# Row 45
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").Value = '159.02'
$ws.Range("E45").Value = '  -2.94%  '

# Row 46
$ws.Range("D46").Value = '192.20'
$ws.Range("E46").Value = '  -1.85%  '

# Row 47
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.04%  '

# Row 48
$ws.Range("D48").Value = '1.36'
$ws.Range("E48").Value = '  +6.15%  '

# Row 49
$ws.Range("D49").Value = '0.790'
$ws.Range("E49").Value = '  -0.80%  '

# Row 50
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '43.39'
$ws.Range("E50").Value = '  +2.26%  '

# Row 51
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").Value = '4.34'
$ws.Range("E51").Value = '  +6.50%  '
